$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.144369602203369
$ws.Range("B1").Value = 2.240702390670776
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.169145822525024
$ws.Range("E1").Value = 1.071166634559631
